$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.035749617779885
$ws.Range("D2").Value = 0.004954460565098628
$ws.Range("E2").Value = 0.600233854077544
$ws.Range("F2").Value = 0.5949479694702262
$ws.Range("G2").Value = 0.00237234260763883
$ws.Range("L2").Value = 0.188535276634255
$ws.Range("M2").Value = 0.2311058523306144
$ws.Range("N2").Value = 1.788967750967686
$ws.Range("O2").Value = 1.920498940323597

$ws.Range("B3").Value = 0.9981439650333925
$ws.Range("D3").Value = 0.004601835448667657
$ws.Range("E3").Value = 0.5881624506335044
$ws.Range("F3").Value = 0.5725732467820848
$ws.Range("G3").Value = 0.002375707973692993
$ws.Range("L3").Value = 0.169926616804176
$ws.Range("M3").Value = 0.2178079834787567
$ws.Range("N3").Value = 1.783896501362975
$ws.Range("O3").Value = 1.857999587130024

$ws.Range("B4").Value = 0.9754366158656183
$ws.Range("D4").Value = 0.00438426586616103
$ws.Range("E4").Value = 0.580783579254792
$ws.Range("F4").Value = 0.5592695690712048
$ws.Range("G4").Value = 0.002377886805572755
$ws.Range("L4").Value = 0.1585149031715645
$ws.Range("M4").Value = 0.2097107768555304
$ws.Range("N4").Value = 1.781496194660306
$ws.Range("O4").Value = 1.8210650989393

$ws.Range("B5").Value = 0.966280173701108
$ws.Range("D5").Value = 0.004295345443146914
$ws.Range("E5").Value = 0.5777853708718794
$ws.Range("F5").Value = 0.5539571654938982
$ws.Range("G5").Value = 0.002378803070049604
$ws.Range("L5").Value = 0.1538683350912748
$ws.Range("M5").Value = 0.2064283532219378
$ws.Range("N5").Value = 1.780698286435651
$ws.Range("O5").Value = 1.806375100828774

$ws.Range("B6").Value = 0.9647656336550483
$ws.Range("D6").Value = 0.004280564822501276
$ws.Range("E6").Value = 0.5772880615380558
$ws.Range("F6").Value = 0.5530816165723564
$ws.Range("G6").Value = 0.002378956931438541
$ws.Range("L6").Value = 0.1530970129886668
$ws.Range("M6").Value = 0.2058843578656919
$ws.Range("N6").Value = 1.780576704774489
$ws.Range("O6").Value = 1.80395761602847

$ws.Range("B7").Value = 0.9753127351802107
$ws.Range("D7").Value = 0.004383067694778475
$ws.Range("E7").Value = 0.580743108388738
$ws.Range("F7").Value = 0.559197483395721
$ws.Range("G7").Value = 0.002377899047685299
$ws.Range("L7").Value = 0.1584522222097888
$ws.Range("M7").Value = 0.209666438790407
$ws.Range("N7").Value = 1.78148470291454
$ws.Range("O7").Value = 1.82086552383052

$ws.Range("B8").Value = 1.022704206908173
$ws.Range("D8").Value = 0.004833098492738941
$ws.Range("E8").Value = 0.5960650399248237
$ws.Range("F8").Value = 0.5871427951886119
$ws.Range("G8").Value = 0.002373479692844072
$ws.Range("L8").Value = 0.1821162372493603
$ws.Range("M8").Value = 0.2265068156022991
$ws.Range("N8").Value = 1.787071595420713
$ws.Range("O8").Value = 1.898649537183246

$ws.Range("B9").Value = 1.118647186434487
$ws.Range("D9").Value = 0.005706979093552889
$ws.Range("E9").Value = 0.6263555380765879
$ws.Range("F9").Value = 0.6454082660964957
$ws.Range("G9").Value = 0.002365701852242994
$ws.Range("L9").Value = 0.2286240124896608
$ws.Range("M9").Value = 0.2600603300521342
$ws.Range("N9").Value = 1.803654578278397
$ws.Range("O9").Value = 2.062671894868316

$ws.Range("B10").Value = 1.190941165194147
$ws.Range("D10").Value = 0.006343488843619127
$ws.Range("E10").Value = 0.648737228185631
$ws.Range("F10").Value = 0.6903572764116319
$ws.Range("G10").Value = 0.002360523523879047
$ws.Range("L10").Value = 0.2628472335487686
$ws.Range("M10").Value = 0.2850269892408051
$ws.Range("N10").Value = 1.819226791524514
$ws.Range("O10").Value = 2.190280014115444

$ws.Range("B11").Value = 1.224215733566723
$ws.Range("D11").Value = 0.006631799806985583
$ws.Range("E11").Value = 0.6589424325347508
$ws.Range("F11").Value = 0.7112772889553298
$ws.Range("G11").Value = 0.002358282974954083
$ws.Range("L11").Value = 0.278426361555006
$ws.Range("M11").Value = 0.2964518255866224
$ws.Range("N11").Value = 1.827038348155043
$ws.Range("O11").Value = 2.249895980919121

$ws.Range("B12").Value = 1.236871039891696
$ws.Range("D12").Value = 0.006740791547720448
$ws.Range("E12").Value = 0.6628098627706365
$ws.Range("F12").Value = 0.7192675108927205
$ws.Range("G12").Value = 0.002357450997944413
$ws.Range("L12").Value = 0.2843271139732053
$ws.Range("M12").Value = 0.3007876167408412
$ws.Range("N12").Value = 1.830100265568959
$ws.Range("O12").Value = 2.272697675539575

$ws.Range("B13").Value = 1.234143062031762
$ws.Range("D13").Value = 0.006717326580915284
$ws.Range("E13").Value = 0.6619768190100928
$ws.Range("F13").Value = 0.7175436322683453
$ws.Range("G13").Value = 0.002357629447869029
$ws.Range("L13").Value = 0.2830562291187562
$ws.Range("M13").Value = 0.2998534098277474
$ws.Range("N13").Value = 1.829436220558009
$ws.Range("O13").Value = 2.267776833739219

$ws.Range("B14").Value = 1.225255796529098
$ws.Range("D14").Value = 0.006640770377963179
$ws.Range("E14").Value = 0.6592605520469945
$ws.Range("F14").Value = 0.7119332792888144
$ws.Range("G14").Value = 0.002358214198129135
$ws.Range("L14").Value = 0.2789117960616068
$ws.Range("M14").Value = 0.29680834561605
$ws.Range("N14").Value = 1.827288177112493
$ws.Range("O14").Value = 2.251767345495921

$ws.Range("B15").Value = 1.219819222188789
$ws.Range("D15").Value = 0.006593853147379036
$ws.Range("E15").Value = 0.6575971294145617
$ws.Range("F15").Value = 0.7085056761250996
$ws.Range("G15").Value = 0.002358574516976271
$ws.Range("L15").Value = 0.2763733688564116
$ws.Range("M15").Value = 0.2949443800245177
$ws.Range("N15").Value = 1.825985939614355
$ws.Range("O15").Value = 2.241990593372407

$ws.Range("B16").Value = 1.188774326668494
$ws.Range("D16").Value = 0.006324621483891235
$ws.Range("E16").Value = 0.6480707307040348
$ws.Range("F16").Value = 0.6889996397316906
$ws.Range("G16").Value = 0.00236067225819947
$ws.Range("L16").Value = 0.2618292945138023
$ws.Range("M16").Value = 0.2842816850065617
$ws.Range("N16").Value = 1.818730863043001
$ws.Range("O16").Value = 2.186415590917363

$ws.Range("B17").Value = 1.16982801342175
$ws.Range("D17").Value = 0.006159133871989297
$ws.Range("E17").Value = 0.6422323248454376
$ws.Range("F17").Value = 0.6771545636614889
$ws.Range("G17").Value = 0.002361988576272705
$ws.Range("L17").Value = 0.2529095458726545
$ws.Range("M17").Value = 0.2777575564787327
$ws.Range("N17").Value = 1.81446592036994
$ws.Range("O17").Value = 2.152724155460646

$ws.Range("B18").Value = 1.158967132498105
$ws.Range("D18").Value = 0.006063833424377663
$ws.Range("E18").Value = 0.6388764828183966
$ws.Range("F18").Value = 0.6703860420997501
$ws.Range("G18").Value = 0.002362756525929403
$ws.Range("L18").Value = 0.2477801768455379
$ws.Range("M18").Value = 0.2740114062544023
$ws.Range("N18").Value = 1.812081377288379
$ws.Range("O18").Value = 2.133493099009286

$ws.Range("B19").Value = 1.155296124499046
$ws.Range("D19").Value = 0.006031546561697354
$ws.Range("E19").Value = 0.6377406536921768
$ws.Range("F19").Value = 0.6681019619177135
$ws.Range("G19").Value = 0.002363018404640494
$ws.Range("L19").Value = 0.2460436468496852
$ws.Range("M19").Value = 0.2727441229752898
$ws.Range("N19").Value = 1.811285807823609
$ws.Range("O19").Value = 2.12700706813996

$ws.Range("B20").Value = 1.171841105115874
$ws.Range("D20").Value = 0.006176762403054425
$ws.Range("E20").Value = 0.6428536031250403
$ws.Range("F20").Value = 0.6784108891647662
$ws.Range("G20").Value = 0.002361847330898348
$ws.Range("L20").Value = 0.2538589633108188
$ws.Range("M20").Value = 0.2784514049870879
$ws.Range("N20").Value = 1.814912842996293
$ws.Range("O20").Value = 2.156295405908565

$ws.Range("B21").Value = 1.22786471687391
$ws.Range("D21").Value = 0.006663261876660442
$ws.Range("E21").Value = 0.6600583093543122
$ws.Range("F21").Value = 0.7135793206187344
$ws.Range("G21").Value = 0.002358041996548979
$ws.Range("L21").Value = 0.2801290841450452
$ws.Range("M21").Value = 0.2977025000931164
$ws.Range("N21").Value = 1.827916297312726
$ws.Range("O21").Value = 2.256463565877652

$ws.Range("B22").Value = 1.264799304281837
$ws.Range("D22").Value = 0.006980133762432672
$ws.Range("E22").Value = 0.6713195079071426
$ws.Range("F22").Value = 0.7369619872945634
$ws.Range("G22").Value = 0.002355650954301602
$ws.Range("L22").Value = 0.2973054275666698
$ws.Range("M22").Value = 0.3103391975468028
$ws.Range("N22").Value = 1.837019766584262
$ws.Range("O22").Value = 2.323249553441997

$ws.Range("B23").Value = 1.245057610395747
$ws.Range("D23").Value = 0.006811114719525335
$ws.Range("E23").Value = 0.6653077974058803
$ws.Range("F23").Value = 0.7244456880648329
$ws.Range("G23").Value = 0.002356918344503319
$ws.Range("L23").Value = 0.2881375141423064
$ws.Range("M23").Value = 0.3035897999320269
$ws.Range("N23").Value = 1.832105978131665
$ws.Range("O23").Value = 2.287483418752174

$ws.Range("B24").Value = 1.170930887844179
$ws.Range("D24").Value = 0.006168793039986298
$ws.Range("E24").Value = 0.6425727208941723
$ws.Range("F24").Value = 0.6778427755620271
$ws.Range("G24").Value = 0.002361911153171934
$ws.Range("L24").Value = 0.2534297356743309
$ws.Range("M24").Value = 0.2781377015482533
$ws.Range("N24").Value = 1.814710579178609
$ws.Range("O24").Value = 2.154680411937875

$ws.Range("B25").Value = 1.092373036584348
$ws.Range("D25").Value = 0.00547152308983101
$ws.Range("E25").Value = 0.6181372506418228
$ws.Range("F25").Value = 0.6292719677446854
$ws.Range("G25").Value = 0.002367711425439345
$ws.Range("L25").Value = 0.2160323292169721
$ws.Range("M25").Value = 0.2509273445481668
$ws.Range("N25").Value = 1.798570372149001
$ws.Range("O25").Value = 2.017059542957782

